# remove some redundant codes
# Appends three more days (9/28, 9/29, 9/30) of 戰隊戰 (clan-battle) log
# entries to 工作表1, mirroring the existing day-block layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dateFmt = 'm"月"d"日"'

# ---- 9/28 (day block: rows 99-111) ----------------------------------
$ws.Range("A99").Value = 43371
$ws.Range("A99").NumberFormat = $dateFmt

$ws.Range("A100").Value = '凱留'
$ws.Range("B100").Value = '亞里莎(UB5)'
$ws.Range("C100").Value = '惠理子(UB2)'
$ws.Range("D100").Value = '日和'
$ws.Range("E100").Value = '純(借五)'
$ws.Range("F100").Value = 255129
$ws.Range("G100").Value = '原一三隊互換與捨棄月月借四星犬好像不一定比較好?'

$ws.Range("A101").Value = 81
$ws.Range("B101").Value = 81
$ws.Range("C101").Value = 81
$ws.Range("D101").Value = 81
$ws.Range("E101").Value = 78
$ws.Range("F101").Value = '龍二'

$ws.Range("A102").Value = 33890
$ws.Range("B102").Value = 67486
$ws.Range("C102").Value = 81793
$ws.Range("D102").Value = 60450
$ws.Range("E102").Value = 11510

$ws.Range("A104").Value = '病弓(UB3)'
$ws.Range("B104").Value = '可可蘿'
$ws.Range("C104").Value = '月月'
$ws.Range("D104").Value = '琉球犬'
$ws.Range("E104").Value = '純(借四)'
$ws.Range("F104").Value = 464926

$ws.Range("A105").Value = 81
$ws.Range("B105").Value = 81
$ws.Range("C105").Value = 81
$ws.Range("D105").Value = 81
$ws.Range("E105").Value = 81
$ws.Range("F105").Value = '龍二'

$ws.Range("A106").Value = 132869
$ws.Range("B106").Value = 43984
$ws.Range("C106").Value = 89846
$ws.Range("D106").Value = 180740
$ws.Range("E106").Value = 17487

$ws.Range("A108").Value = '爆弓(UB2)'
$ws.Range("B108").Value = '美美(UB2)'
$ws.Range("C108").Value = '深月'
$ws.Range("D108").Value = '琉球犬(借四)'
$ws.Range("E108").Value = '純'
$ws.Range("F108").Value = 339827

$ws.Range("A109").Value = 81
$ws.Range("B109").Value = 81
$ws.Range("C109").Value = 81
$ws.Range("D109").Value = 81
$ws.Range("E109").Value = 81
$ws.Range("F109").Value = '龍二'

$ws.Range("A110").Value = 89727
$ws.Range("B110").Value = 68660
$ws.Range("C110").Value = 31827
$ws.Range("D110").Value = 137266
$ws.Range("E110").Value = 12347

$ws.Range("F111").Formula = "=SUM(F100:F109)"

# ---- 9/29 (day block: rows 114-126) ----------------------------------
$ws.Range("A114").Value = 43372
$ws.Range("A114").NumberFormat = $dateFmt

$ws.Range("A115").Value = '凱留'
$ws.Range("B115").Value = '栞'
$ws.Range("C115").Value = '月月(借四)'
$ws.Range("D115").Value = '望'
$ws.Range("E115").Value = '布丁'
$ws.Range("F115").Value = 220766

$ws.Range("A116").Value = 81
$ws.Range("B116").Value = 81
$ws.Range("C116").Value = 81
$ws.Range("D116").Value = 81
$ws.Range("E116").Value = 78
$ws.Range("F116").Value = '豬二'

$ws.Range("A117").Value = 35703
$ws.Range("B117").Value = 90789
$ws.Range("C117").Value = 67119
$ws.Range("D117").Value = 14122
$ws.Range("E117").Value = 13033

$ws.Range("A119").Value = '真步'
$ws.Range("B119").Value = '深月'
$ws.Range("C119").Value = '可可蘿'
$ws.Range("D119").Value = '琉球犬'
$ws.Range("E119").Value = '純(借五)'
$ws.Range("F119").Value = 209698

$ws.Range("A120").Value = 81
$ws.Range("B120").Value = 81
$ws.Range("C120").Value = 81
$ws.Range("D120").Value = 81
$ws.Range("E120").Value = 81
$ws.Range("F120").Value = '豬二'

$ws.Range("A121").Value = 11417
$ws.Range("B121").Value = 29622
$ws.Range("C121").Value = 35400
$ws.Range("D121").Value = 120614
$ws.Range("E121").Value = 12645

$ws.Range("A123").Value = '爆弓(UB2)'
$ws.Range("B123").Value = '美美(UB2)'
$ws.Range("C123").Value = '深月'
$ws.Range("D123").Value = '琉球犬(借四)'
$ws.Range("E123").Value = '純'
$ws.Range("F123").Value = 339827

$ws.Range("A124").Value = 81
$ws.Range("B124").Value = 81
$ws.Range("C124").Value = 81
$ws.Range("D124").Value = 81
$ws.Range("E124").Value = 81
$ws.Range("F124").Value = '龍二'

$ws.Range("A125").Value = 89727
$ws.Range("B125").Value = 68660
$ws.Range("C125").Value = 31827
$ws.Range("D125").Value = 200000
$ws.Range("E125").Value = 12347

$ws.Range("F126").Formula = "=SUM(F115:F124)"

# blank separator row, but keeps the date-number-format stamp
$ws.Range("A128").NumberFormat = $dateFmt

# ---- 9/30 (day block: rows 129-141) ----------------------------------
$ws.Range("A129").Value = 43373
$ws.Range("A129").NumberFormat = $dateFmt

$ws.Range("A130").Value = '栞'
$ws.Range("B130").Value = '惠理子'
$ws.Range("C130").Value = '月月(借四)'
$ws.Range("D130").Value = '望'
$ws.Range("E130").Value = '純'
$ws.Range("F130").Value = 305305

$ws.Range("A131").Value = 82
$ws.Range("B131").Value = 82
$ws.Range("C131").Value = 78
$ws.Range("D131").Value = 82
$ws.Range("E131").Value = 82
$ws.Range("F131").Value = '鹿二'

$ws.Range("A132").Value = 139313
$ws.Range("B132").Value = 101079
$ws.Range("C132").Value = 44387
$ws.Range("D132").Value = 11033
$ws.Range("E132").Value = 9493

$ws.Range("A134").Value = '優衣'
$ws.Range("B134").Value = '可可蘿'
$ws.Range("C134").Value = '琉球犬'
$ws.Range("D134").Value = '月月'
$ws.Range("E134").Value = '純(借五)'
$ws.Range("F134").Value = 303519
$ws.Range("G134").Value = '放栞 15秒犬會死 5秒366647'

$ws.Range("A135").Value = 82
$ws.Range("B135").Value = 82
$ws.Range("C135").Value = 82
$ws.Range("D135").Value = 82
$ws.Range("E135").Value = 80
$ws.Range("F135").Value = '鹿二'
$ws.Range("G135").Value = '原隊伍 4秒 298330'

$ws.Range("A136").Value = 24206
$ws.Range("B136").Value = 51506
$ws.Range("C136").Value = 145790
$ws.Range("D136").Value = 65251
$ws.Range("E136").Value = 16766

$ws.Range("A138").Value = '爆弓'
$ws.Range("B138").Value = '亞里莎'
$ws.Range("C138").Value = '深月'
$ws.Range("D138").Value = '優花梨'
$ws.Range("E138").Value = '純(借四)'
$ws.Range("F138").Value = 226616

$ws.Range("A139").Value = 82
$ws.Range("B139").Value = 82
$ws.Range("C139").Value = 82
$ws.Range("D139").Value = 82
$ws.Range("E139").Value = 79
$ws.Range("F139").Value = '鹿二'

$ws.Range("A140").Value = 100794
$ws.Range("B140").Value = 90952
$ws.Range("C140").Value = 19049
$ws.Range("D140").Value = 7973
$ws.Range("E140").Value = 7848

$ws.Range("F141").Formula = "=SUM(F130:F139)"

# trailing blank row stamped with the same date number format
$ws.Range("A144").NumberFormat = $dateFmt

# ---- view state: scroll down to the newly-added rows -----------------
$ws.Activate() | Out-Null
$ws.Range("G124").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 112
$excel.ActiveWindow.ScrollColumn = 1
